$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)

$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

$ws.Range("A2").Value = 91
$ws.Range("B2").Value = "富邦人壽"
$ws.Range("C2").Value = "身壽險丙型"
$ws.Range("D2").Value = "孫大千"
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
$ws.Range("G2").Value = "2011-11-21"
$ws.Range("H2").Value = "孫大千"
$ws.Range("I2").Value = 919
$ws.Range("J2").Value = "tmpc6841"
$ws.Range("K2").Value = 91

$ws.Range("A3").Value = 93
$ws.Range("B3").Value = "富邦人壽"
$ws.Range("C3").Value = "富邦人壽重大疾病終身保險"
$ws.Range("D3").Value = "孫大千"
$ws.Range("E3").Value = "insurance"
$ws.Range("F3").Value = "normal"
$ws.Range("G3").Value = "2011-11-21"
$ws.Range("H3").Value = "孫大千"
$ws.Range("I3").Value = 919
$ws.Range("J3").Value = "tmpc6841"
$ws.Range("K3").Value = 93

$ws.Range("A4").Value = 96
$ws.Range("B4").Value = "富邦人壽"
$ws.Range("C4").Value = "富邦人壽分红養老保險"
$ws.Range("D4").Value = "孫大千"
$ws.Range("E4").Value = "insurance"
$ws.Range("F4").Value = "normal"
$ws.Range("G4").Value = "2011-11-21"
$ws.Range("H4").Value = "孫大千"
$ws.Range("I4").Value = 919
$ws.Range("J4").Value = "tmpc6841"
$ws.Range("K4").Value = 96

$ws.Range("A5").Value = 99
$ws.Range("B5").Value = "富邦人壽"
$ws.Range("C5").Value = "富邦人壽雙福還本分紅終身"
$ws.Range("D5").Value = "陳端梅"
$ws.Range("E5").Value = "insurance"
$ws.Range("F5").Value = "normal"
$ws.Range("G5").Value = "2011-11-21"
$ws.Range("H5").Value = "孫大千"
$ws.Range("I5").Value = 919
$ws.Range("J5").Value = "tmpc6841"
$ws.Range("K5").Value = 99

$ws.Range("A6").Value = 103
$ws.Range("B6").Value = "保德信人壽"
$ws.Range("C6").Value = "保德信教育終身壽險"
$ws.Range("D6").Value = "陳端梅"
$ws.Range("E6").Value = "insurance"
$ws.Range("F6").Value = "normal"
$ws.Range("G6").Value = "2011-11-21"
$ws.Range("H6").Value = "孫大千"
$ws.Range("I6").Value = 919
$ws.Range("J6").Value = "tmpc6841"
$ws.Range("K6").Value = 103

$ws.Range("A7").Value = 104
$ws.Range("B7").Value = "保德信人壽"
$ws.Range("C7").Value = "保德信特定傷病終身健康保險"
$ws.Range("D7").Value = "陳端梅"
$ws.Range("E7").Value = "insurance"
$ws.Range("F7").Value = "normal"
$ws.Range("G7").Value = "2011-11-21"
$ws.Range("H7").Value = "孫大千"
$ws.Range("I7").Value = 919
$ws.Range("J7").Value = "tmpc6841"
$ws.Range("K7").Value = 104

$ws.Range("A8").Value = 105
$ws.Range("B8").Value = "保德信人壽"
$ws.Range("C8").Value = "保德信癌症終身健康保險"
$ws.Range("D8").Value = "陳端梅"
$ws.Range("E8").Value = "insurance"
$ws.Range("F8").Value = "normal"
$ws.Range("G8").Value = "2011-11-21"
$ws.Range("H8").Value = "孫大千"
$ws.Range("I8").Value = 919
$ws.Range("J8").Value = "tmpc6841"
$ws.Range("K8").Value = 105

$ws.Range("A9").Value = 106
$ws.Range("B9").Value = "富#人壽"
$ws.Range("C9").Value = "富邦人壽全福還本分紅终身壽險"
$ws.Range("D9").Value = "陳端梅"
$ws.Range("E9").Value = "insurance"
$ws.Range("F9").Value = "normal"
$ws.Range("G9").Value = "2011-11-21"
$ws.Range("H9").Value = "孫大千"
$ws.Range("I9").Value = 919
$ws.Range("J9").Value = "tmpc6841"
$ws.Range("K9").Value = 106
